# Update "paises" (countries) COVID data sheet + provincias Spain
# per commit "Update countries & provincias Spain"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "last updated" timestamp banner in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 12:25"

# --- 2. Refresh per-country case numbers (Casos totales, Nuevos casos,
#         Casos activos, Recuperados, Casos criticos, Muertes) ---

# Row 15 - Iran
$ws.Range("B15").Value = 388810
$ws.Range("C15").Value = 2152
$ws.Range("D15").Value = 335572
$ws.Range("E15").Value = 30828
$ws.Range("G15").Value = 117
$ws.Range("H15").Value = 22410

# Row 17 - Banglades
$ws.Range("B17").Value = 327359
$ws.Range("C17").Value = 2202
$ws.Range("D17").Value = 224573
$ws.Range("E17").Value = 98270
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = 4516

# Row 24 - Alemania
$ws.Range("B24").Value = 251744
$ws.Range("C24").Value = 20
$ws.Range("E24").Value = 15343

# Row 37 - Rumania
$ws.Range("B37").Value = 95897
$ws.Range("C37").Value = 883
$ws.Range("D37").Value = 40454
$ws.Range("E37").Value = 51517
$ws.Range("G37").Value = 33
$ws.Range("H37").Value = 3926

# Row 40 - Oman
$ws.Range("B40").Value = 87328
$ws.Range("C40").Value = 256
$ws.Range("D40").Value = 82805
$ws.Range("E40").Value = 3789
$ws.Range("G40").Value = 6
$ws.Range("H40").Value = 734

# Row 45 - Emiratos Arabes Unidos
$ws.Range("B45").Value = 74454
$ws.Range("C45").Value = 470
$ws.Range("D45").Value = 66533
$ws.Range("E45").Value = 7531
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 390

# Row 96 - Malasia
$ws.Range("B96").Value = 9459
$ws.Range("C96").Value = 62
$ws.Range("D96").Value = 9124
$ws.Range("E96").Value = 207

# Row 104 - Finlandia
$ws.Range("B104").Value = 8327
$ws.Range("C104").Value = 36
$ws.Range("E104").Value = 641

# --- 3. Eslovenia overtakes Siria in the ranking (sorted by Casos totales
#         descending), so the two rows swap places; Siria's own figures are
#         unchanged, Eslovenia's are refreshed ---

# Row 129 now holds Eslovenia (previously Siria) with updated figures
$ws.Range("A129").Value = "Eslovenia"
$ws.Range("B129").Value = 3190
$ws.Range("C129").Value = 25
$ws.Range("D129").Value = 2530
$ws.Range("E129").Value = 525
$ws.Range("H129").Value = 135

# Row 130 now holds Siria (previously Eslovenia) with its prior figures
$ws.Range("A130").Value = "Siria"
$ws.Range("B130").Value = 3171
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 730
$ws.Range("E130").Value = 2307
$ws.Range("H130").Value = 134
